$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Row 2 (Beta) ---
$ws.Range("C2").Value = 16.27395536621195
$ws.Range("E2").Value = 0.006275575909959944
$ws.Range("F2").Value = 8.324375026313408
$ws.Range("G2").Value = 7.920849506013032
$ws.Range("H2").Value = 8.777728043786844
$ws.Range("I2").Value = 0.100050592873538
$ws.Range("J2").Value = 0.09803150746425697
$ws.Range("K2").Value = 0.1024391192949508
$ws.Range("L2").Value = 0.01317952011340569
$ws.Range("M2").Value = 0.01283356493344551
$ws.Range("N2").Value = 0.01358275867402604

# --- Update Row 3 (Gamma) ---
$ws.Range("C3").Value = 0.3589057182506037
$ws.Range("D3").Value = 0.3038166771491592
$ws.Range("E3").Value = 0.3557873748505794
$ws.Range("F3").Value = 0.2882425796302744
$ws.Range("G3").Value = 0.1986911320191278
$ws.Range("H3").Value = 0.3723267416218894
$ws.Range("I3").Value = 0.2669917055821284
$ws.Range("J3").Value = 0.1850185589117773
$ws.Range("K3").Value = 0.3439887162973743
$ws.Range("L3").Value = 0.2938553798918609
$ws.Range("M3").Value = 0.198306971834185
$ws.Range("N3").Value = 0.3834988144254472

# --- Add new Row 4 (Beta + Gamma) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 16.63286108446255
$ws.Range("D4").Value = 0.3068349094013313
$ws.Range("E4").Value = 0.3620629507605393
$ws.Range("F4").Value = 8.612617605943683
$ws.Range("G4").Value = 8.119540638032159
$ws.Range("H4").Value = 9.150054785408734
$ws.Range("I4").Value = 0.3670422984556664
$ws.Range("J4").Value = 0.2830500663760343
$ws.Range("K4").Value = 0.4464278355923251
$ws.Range("L4").Value = 0.3070349000052666
$ws.Range("M4").Value = 0.2111405367676305
$ws.Range("N4").Value = 0.3970815730994733

# Copy the number/style formatting of A2 (which matches A3) onto the new A4 cell
# so it picks up the same cell style index as the existing first-column cells.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
